$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
$c = $tcs.Colors(3)
Write-Output "before:" $c.RGB
$c.RGB = 0x123456
Write-Output "after:" $tcs.Colors(3).RGB
